$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.077.67"
$ws.Range("E2").Value = "  -1.89%  "

$ws.Range("D3").Value = "1.832.57"
$ws.Range("E3").Value = "  -0.92%  "

$ws.Range("E4").Value = "  -0.08%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "325.18"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -3.12%  "

$ws.Range("E6").Value = "  -0.09%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4617"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.98%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3857"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.29%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07845"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.65%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.9600"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.27%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "21.94"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.63%  "

$ws.Range("D12").Value = "1.883.13"
$ws.Range("E12").Value = "  +1.63%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.670"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -2.91%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.881"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.80%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.06846"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.05%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "88.13"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.54%  "

$ws.Range("E17").Value = "  -0.06%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000009911"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.20%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "16.65"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.76%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").Value = "28.099.22"
$ws.Range("E21").Value = "  -1.84%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.287"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.00%  "

$ws.Range("E23").Value = "  -2.67%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.085"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -3.82%  "

$ws.Range("D25").Value = "2.099.99"
$ws.Range("E25").Value = "  +1.11%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "154.43"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.47%  "

$ws.Range("E27").Value = "  -1.50%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "5.748"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -5.87%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.966"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -3.17%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "118.59"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.84%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.9420"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -3.40%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.09225"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.43%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.263"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.74%  "

$ws.Range("E34").Value = "  -2.14%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "3.326"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -4.53%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.05820"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -5.47%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.02108"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -4.12%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.135"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.46%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "7.716"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.46%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.5589"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.36%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "9.888"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.72%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1759"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.10%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.07314"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.91%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "11.63"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.05%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.5263"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.37%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.141"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -8.46%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.105"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -11.13%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.832"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -3.90%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "113.40"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.40%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.10%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.021"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.03%  "
